# Baseball Catcher trial runs: record the results of "Run 1" and finish
# filling in "Run 2", then point the workbook back at "Run 1" as the
# active sheet with the new F column legend cell selected.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Run 2: mark three more rims (column B) for rows 2-4.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Run 2")
$ws2.Range("B2").Value = 1
$ws2.Range("B3").Value = 1
$ws2.Range("B4").Value = 1
$ws2.Range("B6").Select()

# ---------------------------------------------------------------
# Run 1: insert a new summary row above the table, pushing the
# existing header + data down by one row, then log one more catch.
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Run 1")
$ws1.Rows.Item(1).Insert()

# Summary formulas in the new row 1: total points from catches (x4),
# rims (x2), a flat 0 for misses, and the grand total.
$ws1.Range("A1").Formula = "=SUM(A3:A23)*4"
$ws1.Range("B1").Formula = "=SUM(B3:B23)*2"
$ws1.Range("C1").Value = 0
$ws1.Range("D1").Formula = "=SUM(A1:B1)"

# One more catch thrown during the session.
$ws1.Range("A23").Value = 1

# Legend in column F summarizing the tallies for this run.
$ws1.Range("F2").Value = "11 catches"
$ws1.Range("F3").Value = "8 rims"
$ws1.Range("F4").Value = "1 miss"
$ws1.Range("F5").Value = "1 uncatchable"

$ws1.Activate()
$ws1.Range("F2").Select()
